$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.05952403694391251
$ws.Cells.Item(2, 2).Value = 0.9840734601020813
$ws.Cells.Item(2, 3).Value = 0.03294847160577774
$ws.Cells.Item(2, 4).Value = 0.995579719543457
$ws.Cells.Item(3, 1).Value = 0.009315329603850842
$ws.Cells.Item(3, 2).Value = 0.9985074400901794
$ws.Cells.Item(3, 3).Value = 0.02383473142981529
$ws.Cells.Item(3, 4).Value = 0.995579719543457
$ws.Cells.Item(4, 1).Value = 0.004900739993900061
$ws.Cells.Item(4, 2).Value = 0.9988410472869873
$ws.Cells.Item(4, 3).Value = 0.0116785941645503
$ws.Cells.Item(4, 4).Value = 0.9986399412155151
$ws.Cells.Item(5, 1).Value = 0.002231738762930036
$ws.Cells.Item(5, 2).Value = 0.9993854165077209
$ws.Cells.Item(5, 3).Value = 0.009853780269622803
$ws.Cells.Item(5, 4).Value = 0.9986399412155151
$ws.Cells.Item(6, 1).Value = 0.00176722917240113
$ws.Cells.Item(6, 2).Value = 0.9994907975196838
$ws.Cells.Item(6, 3).Value = 0.0007578931981697679
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(7, 1).Value = 0.001555059570819139
$ws.Cells.Item(7, 2).Value = 0.9996312260627747
$ws.Cells.Item(7, 3).Value = 0.0004766620113514364
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(8, 1).Value = 0.001264742226339877
$ws.Cells.Item(8, 2).Value = 0.9996312260627747
$ws.Cells.Item(8, 3).Value = 0.0006716083735227585
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(9, 1).Value = 0.001739514176733792
$ws.Cells.Item(9, 2).Value = 0.9996137022972107
$ws.Cells.Item(9, 3).Value = 0.001007005921564996
$ws.Cells.Item(9, 4).Value = 0.9993199706077576
$ws.Cells.Item(10, 1).Value = 0.0008868636796250939
$ws.Cells.Item(10, 2).Value = 0.9998244047164917
$ws.Cells.Item(10, 3).Value = 0.0005952870124019682
$ws.Cells.Item(10, 4).Value = 0.9996599555015564
$ws.Cells.Item(11, 1).Value = 0.0006371597992256284
$ws.Cells.Item(11, 2).Value = 0.9998419880867004
$ws.Cells.Item(11, 3).Value = 0.002753538312390447
$ws.Cells.Item(11, 4).Value = 0.9993199706077576
$ws.Cells.Item(12, 1).Value = 0.001135466271080077
$ws.Cells.Item(12, 2).Value = 0.9997190237045288
$ws.Cells.Item(12, 3).Value = 0.0001414724247297272
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(13, 1).Value = 0.00044366589281708
$ws.Cells.Item(13, 2).Value = 0.9998946189880371
$ws.Cells.Item(13, 3).Value = 0.0000476097411592491
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(14, 1).Value = 0.0003509992093313485
$ws.Cells.Item(14, 2).Value = 0.9999122023582458
$ws.Cells.Item(14, 3).Value = 0.00003694478436955251
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(15, 1).Value = 0.000989259104244411
$ws.Cells.Item(15, 2).Value = 0.9996839165687561
$ws.Cells.Item(15, 3).Value = 0.00006552002741955221
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(16, 1).Value = 0.0001409975229762495
$ws.Cells.Item(16, 2).Value = 0.9999824166297913
$ws.Cells.Item(16, 3).Value = 0.00003014301182702184
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(17, 1).Value = 0.0002642763138283044
$ws.Cells.Item(17, 2).Value = 0.9999297857284546
$ws.Cells.Item(17, 3).Value = 0.00002164144098060206
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(18, 1).Value = 0.0004353784024715424
$ws.Cells.Item(18, 2).Value = 0.9998946189880371
$ws.Cells.Item(18, 3).Value = 0.00003993677455582656
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(19, 1).Value = 0.0006414831732399762
$ws.Cells.Item(19, 2).Value = 0.9998770952224731
$ws.Cells.Item(19, 3).Value = 0.00001244471786776558
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(20, 1).Value = 0.0004955293843522668
$ws.Cells.Item(20, 2).Value = 0.9998770952224731
$ws.Cells.Item(20, 3).Value = 0.000001618578153284034
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(21, 1).Value = 0.0003335158980917186
$ws.Cells.Item(21, 2).Value = 0.9998946189880371
$ws.Cells.Item(21, 3).Value = 0.000008271902515843976
$ws.Cells.Item(21, 4).Value = 1
$ws.Cells.Item(22, 1).Value = 0.0005455230129882693
$ws.Cells.Item(22, 2).Value = 0.9998770952224731
$ws.Cells.Item(22, 3).Value = 0.000001277367346119718
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(23, 1).Value = 0.0001853463618317619
$ws.Cells.Item(23, 2).Value = 0.9999824166297913
$ws.Cells.Item(23, 3).Value = 0.000001074944748324924
$ws.Cells.Item(23, 4).Value = 1
$ws.Cells.Item(24, 1).Value = 0.0008273056009784341
$ws.Cells.Item(24, 2).Value = 0.9998595118522644
$ws.Cells.Item(24, 3).Value = 0.0001078614877769724
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(25, 1).Value = 0.000437478709500283
$ws.Cells.Item(25, 2).Value = 0.9998419880867004
$ws.Cells.Item(25, 3).Value = 0.00001591868931427598
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(26, 1).Value = 0.0002174984983867034
$ws.Cells.Item(26, 2).Value = 0.9999473094940186
$ws.Cells.Item(26, 3).Value = 0.000002495477019692771
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(27, 1).Value = 0.00004839718894800171
$ws.Cells.Item(27, 2).Value = 0.9999824166297913
$ws.Cells.Item(27, 3).Value = 0.0000003478202472706471
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(28, 1).Value = 0.0005956218228675425
$ws.Cells.Item(28, 2).Value = 0.9998946189880371
$ws.Cells.Item(28, 3).Value = 0.000002552928435761714
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(29, 1).Value = 0.0005976655520498753
$ws.Cells.Item(29, 2).Value = 0.9998770952224731
$ws.Cells.Item(29, 3).Value = 0.00001937990600708872
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(30, 1).Value = 0.0002494436921551824
$ws.Cells.Item(30, 2).Value = 0.9999648928642273
$ws.Cells.Item(30, 3).Value = 0.000009739025699673221
$ws.Cells.Item(30, 4).Value = 1
$ws.Cells.Item(31, 1).Value = 0.0001364455092698336
$ws.Cells.Item(31, 2).Value = 0.9999824166297913
$ws.Cells.Item(31, 3).Value = 0.000003331159405206563
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(32, 1).Value = 0.0003935116110369563
$ws.Cells.Item(32, 2).Value = 0.9999473094940186
$ws.Cells.Item(32, 3).Value = 0.000001750265823829977
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(33, 1).Value = 0.0003499925951473415
$ws.Cells.Item(33, 2).Value = 0.9999648928642273
$ws.Cells.Item(33, 3).Value = 0.000001911290610223659
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(34, 1).Value = 0.0001907898404169828
$ws.Cells.Item(34, 2).Value = 0.9999824166297913
$ws.Cells.Item(34, 3).Value = 0.00001548062755318824
$ws.Cells.Item(34, 4).Value = 1
$ws.Cells.Item(35, 1).Value = 0.00003466897032922134
$ws.Cells.Item(35, 2).Value = 1
$ws.Cells.Item(35, 3).Value = 0.000001639785637053137
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(36, 1).Value = 0.0004706543695647269
$ws.Cells.Item(36, 2).Value = 0.9999122023582458
$ws.Cells.Item(36, 3).Value = 0.003496474819257855
$ws.Cells.Item(36, 4).Value = 0.998979926109314
$ws.Cells.Item(37, 1).Value = 0.0009024746832437813
$ws.Cells.Item(37, 2).Value = 0.9998946189880371
$ws.Cells.Item(37, 3).Value = 0.0000009133049161391682
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(38, 1).Value = 0.00004606090442393906
$ws.Cells.Item(38, 2).Value = 0.9999824166297913
$ws.Cells.Item(38, 3).Value = 0.0000003700332911193982
$ws.Cells.Item(38, 4).Value = 1
$ws.Cells.Item(39, 1).Value = 0.0003337124362587929
$ws.Cells.Item(39, 2).Value = 0.9998946189880371
$ws.Cells.Item(39, 3).Value = 0.000002365155069128377
$ws.Cells.Item(39, 4).Value = 1
$ws.Cells.Item(40, 1).Value = 0.00003623387237894349
$ws.Cells.Item(40, 2).Value = 1
$ws.Cells.Item(40, 3).Value = 0.00000006201086932833277
$ws.Cells.Item(40, 4).Value = 1
$ws.Cells.Item(41, 1).Value = 0.00001197950768982992
$ws.Cells.Item(41, 2).Value = 1
$ws.Cells.Item(41, 3).Value = 0.00000008255703676240955
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(42, 1).Value = 0.0003608142142184079
$ws.Cells.Item(42, 2).Value = 0.9999473094940186
$ws.Cells.Item(42, 3).Value = 0.00000004527326780134899
$ws.Cells.Item(42, 4).Value = 1
$ws.Cells.Item(43, 1).Value = 0.0003064493939746171
$ws.Cells.Item(43, 2).Value = 0.9998770952224731
$ws.Cells.Item(43, 3).Value = 0.00000002739975535348549
$ws.Cells.Item(43, 4).Value = 1
$ws.Cells.Item(44, 1).Value = 0.0002186317869927734
$ws.Cells.Item(44, 2).Value = 0.9998946189880371
$ws.Cells.Item(44, 3).Value = 0.0000008580321377849032
$ws.Cells.Item(44, 4).Value = 1
$ws.Cells.Item(45, 1).Value = 0.000007329012532863999
$ws.Cells.Item(45, 2).Value = 1
$ws.Cells.Item(45, 3).Value = 0.0000009789381465452607
$ws.Cells.Item(45, 4).Value = 1
$ws.Cells.Item(46, 1).Value = 0.0009627611725591123
$ws.Cells.Item(46, 2).Value = 0.9999122023582458
$ws.Cells.Item(46, 3).Value = 0.0000007829431751815719
$ws.Cells.Item(46, 4).Value = 1
$ws.Cells.Item(47, 1).Value = 0.000201188595383428
$ws.Cells.Item(47, 2).Value = 0.9999297857284546
$ws.Cells.Item(47, 3).Value = 0.0000004813603027287172
$ws.Cells.Item(47, 4).Value = 1
$ws.Cells.Item(48, 1).Value = 0.0006652078009210527
$ws.Cells.Item(48, 2).Value = 0.9998946189880371
$ws.Cells.Item(48, 3).Value = 0.000001048366016220825
$ws.Cells.Item(48, 4).Value = 1
$ws.Cells.Item(49, 1).Value = 0.00004825552969123237
$ws.Cells.Item(49, 2).Value = 0.9999824166297913
$ws.Cells.Item(49, 3).Value = 0.0000003782112116823555
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(50, 1).Value = 0.0001166382280644029
$ws.Cells.Item(50, 2).Value = 0.9999648928642273
$ws.Cells.Item(50, 3).Value = 0.00000002140106580839074
$ws.Cells.Item(50, 4).Value = 1
$ws.Cells.Item(51, 1).Value = 0.0004509545397013426
$ws.Cells.Item(51, 2).Value = 0.9998946189880371
$ws.Cells.Item(51, 3).Value = 0.0000006281791229412192
$ws.Cells.Item(51, 4).Value = 1

Write-Host "Updated training history values"
